{"js": "// Update the worksheet date header and the 25 \"a\u00f7b=\" division prompts\n// inside the table, per the commit diff. Every <w:t> run in the document\n// (title paragraph + one run per populated table cell) gets replaced by\n// its corresponding new value, matched positionally in document order\n// against the OLD text (all old values are unique in this document).\nconst replacements = [\n  [\"2023-11-18 Saturday\", \"2023-11-19 Sunday\"],\n  [\"84\u00f75=\", \"69\u00f78=\"],\n  [\"90\u00f76=\", \"65\u00f72=\"],\n  [\"76\u00f76=\", \"43\u00f76=\"],\n  [\"28\u00f79=\", \"29\u00f76=\"],\n  [\"65\u00f74=\", \"98\u00f76=\"],\n  [\"29\u00f75=\", \"37\u00f76=\"],\n  [\"77\u00f77=\", \"22\u00f76=\"],\n  [\"27\u00f79=\", \"13\u00f76=\"],\n  [\"45\u00f79=\", \"18\u00f78=\"],\n  [\"14\u00f73=\", \"52\u00f74=\"],\n  [\"77\u00f76=\", \"42\u00f79=\"],\n  [\"17\u00f73=\", \"60\u00f75=\"],\n  [\"93\u00f79=\", \"20\u00f73=\"],\n  [\"89\u00f72=\", \"47\u00f79=\"],\n  [\"18\u00f79=\", \"47\u00f76=\"],\n  [\"57\u00f79=\", \"58\u00f76=\"],\n  [\"58\u00f76=\", \"24\u00f79=\"],\n  [\"81\u00f79=\", \"69\u00f73=\"],\n  [\"51\u00f78=\", \"34\u00f72=\"],\n  [\"38\u00f73=\", \"62\u00f74=\"],\n  [\"39\u00f77=\", \"25\u00f79=\"],\n  [\"49\u00f72=\", \"74\u00f77=\"],\n  [\"96\u00f74=\", \"44\u00f76=\"],\n  [\"76\u00f74=\", \"40\u00f76=\"],\n  [\"80\u00f76=\", \"21\u00f77=\"],\n];\n\n// Build a lookup from old text -> queue of new texts (handles the\n// (unlikely here) case of repeated old values by consuming them in\n// document order).\nconst queues = new Map();\nfor (const [oldText, newText] of replacements) {\n  if (!queues.has(oldText)) queues.set(oldText, []);\n  queues.get(oldText).push(newText);\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  const current = paragraph.text;\n  const queue = queues.get(current);\n  if (queue && queue.length) {\n    const newText = queue.shift();\n    paragraph.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date header and the 25 \"a\u00f7b=\" division prompts inside the\n# table, per the commit diff. Cells are addressed positionally\n# (Table.Cell(row, col)) rather than by text search-and-replace, because\n# a couple of the new values collide with OTHER cells' old values (e.g.\n# row 13 col 1 \"57\u00f79=\" -> \"58\u00f76=\" while row 13 col 2 \"58\u00f76=\" -> \"24\u00f79=\");\n# a naive global Find/Replace would clobber the second cell once the\n# first replacement recreates its old text. Positional addressing makes\n# every write independent of current cell contents.\n\n$d = $word.ActiveDocument\n\n# Title paragraph (date line above the table).\n$d.Paragraphs.Item(1).Range.Text = \"2023-11-19 Sunday\"\n\n$table = $d.Tables.Item(1)\n\n# row -> ordered list of new values for columns 1..5\n$rowValues = @{\n    1  = @(\"69\u00f78=\", \"65\u00f72=\", \"43\u00f76=\", \"29\u00f76=\", \"98\u00f76=\")\n    5  = @(\"37\u00f76=\", \"22\u00f76=\", \"13\u00f76=\", \"18\u00f78=\", \"52\u00f74=\")\n    9  = @(\"42\u00f79=\", \"60\u00f75=\", \"20\u00f73=\", \"47\u00f79=\", \"47\u00f76=\")\n    13 = @(\"58\u00f76=\", \"24\u00f79=\", \"69\u00f73=\", \"34\u00f72=\", \"62\u00f74=\")\n    17 = @(\"25\u00f79=\", \"74\u00f77=\", \"44\u00f76=\", \"40\u00f76=\", \"21\u00f77=\")\n}\n\nforeach ($row in $rowValues.Keys) {\n    $values = $rowValues[$row]\n    for ($col = 1; $col -le $values.Count; $col++) {\n        $table.Cell($row, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
